$wb = $excel.ActiveWorkbook

# "想去人数" (want-to-go count) updates in column F, for sheets
# "展览" (sheet1) and "全部类型" (sheet4). Both sheets carry identical
# rows; only F14 differs between the two sheets in the final value.

$updates = @(
    @{ Row = 6;  Sheet1 = 189; Sheet4 = 189 },
    @{ Row = 8;  Sheet1 = 46;  Sheet4 = 46 },
    @{ Row = 11; Sheet1 = 44;  Sheet4 = 44 },
    @{ Row = 13; Sheet1 = 94;  Sheet4 = 94 },
    @{ Row = 14; Sheet1 = 1797; Sheet4 = 1798 },
    @{ Row = 16; Sheet1 = 491; Sheet4 = 491 },
    @{ Row = 17; Sheet1 = 456; Sheet4 = 456 },
    @{ Row = 19; Sheet1 = 70;  Sheet4 = 70 },
    @{ Row = 22; Sheet1 = 1440; Sheet4 = 1440 },
    @{ Row = 23; Sheet1 = 3367; Sheet4 = 3367 },
    @{ Row = 25; Sheet1 = 57;  Sheet4 = 57 },
    @{ Row = 27; Sheet1 = 1090; Sheet4 = 1090 },
    @{ Row = 28; Sheet1 = 81;  Sheet4 = 81 },
    @{ Row = 29; Sheet1 = 1780; Sheet4 = 1780 },
    @{ Row = 30; Sheet1 = 566; Sheet4 = 566 },
    @{ Row = 31; Sheet1 = 457; Sheet4 = 457 },
    @{ Row = 34; Sheet1 = 401; Sheet4 = 401 },
    @{ Row = 36; Sheet1 = 639; Sheet4 = 639 },
    @{ Row = 38; Sheet1 = 35;  Sheet4 = 35 }
)

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

foreach ($u in $updates) {
    $ws1.Range("F" + $u.Row).Value = $u.Sheet1
    $ws4.Range("F" + $u.Row).Value = $u.Sheet4
}
